# Feat: Add Delay Function
#
# Appends a new "DELAY" command row to the RedfishCommands table:
#   A12 = "DELAY"  (Method column)
#   B12 = 10        (Endpoint column, used here as a delay-in-seconds value)
# and leaves the new B12 cell selected, matching the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "DELAY"
$ws.Range("B12").Value = 10

$ws.Range("B12").Select()
